# "write a review testcases for Osprey"
# Adds a new "Review" worksheet (placed after "Minicart") containing the
# review test-data table, matching the author's target workbook layout.
# Cell writes are ordered to reproduce the exact shared-string allocation
# sequence the author's Excel session produced.

$wb = $excel.ActiveWorkbook

# --- add the new sheet at the very end of the tab strip -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Review"

# --- header row (row 1) — columns that reuse strings already present
#     elsewhere in the workbook (carried over from the "Minicart" sheet). ---
$ws.Range("A1").Value = "DataSet"
$ws.Range("B1").Value = "UserName"
$ws.Range("C1").Value = "Prod UserName"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "Password"
$ws.Range("F1").Value = "Confirm Password"
$ws.Range("G1").Value = "FirstName"
$ws.Range("H1").Value = "LastName"
$ws.Range("I1").Value = "Street"
$ws.Range("J1").Value = "City"
$ws.Range("K1").Value = "Region"
$ws.Range("L1").Value = "postcode"
$ws.Range("M1").Value = "phone"
$ws.Range("N1").Value = "Products"
$ws.Range("O1").Value = "Color"
$ws.Range("P1").Value = "Size"
$ws.Range("R1").Value = "Colorproduct"
$ws.Range("S1").Value = "Quantity"
$ws.Range("T1").Value = "methods"
$ws.Range("U1").Value = "cardNumber"
$ws.Range("V1").Value = "ExpMonthYear"
$ws.Range("W1").Value = "cvv"

# --- data rows carried over from the existing DataSet table ---------------
$ws.Range("A2").Value = "Account"
$ws.Range("B2").Value = "testersemail.278@gmail.com"
$ws.Range("D2").Value = "testersemail.278@gmail.com"
$ws.Range("E2").Value = "Testers@278"
$ws.Range("F2").Value = "Testers@278"
$ws.Range("G2").Value = "QA"
$ws.Range("H2").Value = "TEST"
$ws.Range("I2").Value = "6 Walnut Valley Dr"
$ws.Range("J2").Value = "Little Rock"
$ws.Range("K2").Value = "Arkansas"
$ws.Range("L2").Value = 72211
$ws.Range("M2").Value = 9898989898

$ws.Range("A3").Value = "Product"
$ws.Range("N3").Value = "POCO® CARRYING CASE"
$ws.Range("O3").Value = "Black"
$ws.Range("P3").Value = "S/M"
$ws.Range("R3").Value = "AETHER™ 55"
$ws.Range("S3").Value = 1

# --- new review test-data rows --------------------------------------------
$ws.Range("A4").Value = "review"
$ws.Range("B4").Value = "qatesting"
$ws.Range("D4").Value = "qatesting.lotuswave@gmail.com"
$ws.Range("G4").Value = "testing"
$ws.Range("N4").Value = "POCO® CARRYING CASE"
$ws.Range("X4").Value = 2
$ws.Range("Y4").Value = "Average Product"
$ws.Range("Z4").Value = "Prodcut is to expensive"

$ws.Range("A5").Value = "Ask a question"
$ws.Range("D5").Value = "qatesting.lotuswave@gmail.com"
$ws.Range("G5").Value = "TestingQA"
$ws.Range("Q5").Value = "Product Qunatity"

# --- new columns added to the header row for the review fields ------------
$ws.Range("X1").Value = "score"
$ws.Range("Y1").Value = "title"
$ws.Range("Z1").Value = "Review"
$ws.Range("Q1").Value = "CommetsOsprey"

# --- final review row ------------------------------------------------------
$ws.Range("A6").Value = "Review product"
$ws.Range("N6").Value = "POCO® CARRYING CASE"
$ws.Range("S6").Value = "'1"

# --- header row highlight (matches the yellow "DataSet" header style used
#     by every other sheet in the workbook) ---------------------------------
$ws.Range("A1:Z1").Interior.Color = 65535

# --- column widths (best-fit, mirroring the author's manual autosize) -----
foreach ($col in @("B", "C", "D", "G", "I", "X", "Y", "Z")) {
  $ws.Range($col + "1").EntireColumn.AutoFit() | Out-Null
}

# --- selection / active cell -----------------------------------------------
$ws.Range("N4:N5").Select()
